$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text values (e.g. "573.46", "0.0670").
# Force text number format so Excel does not silently convert them to real numbers
# (which would drop significant trailing zeros / change formatting).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "66.759.77"
$ws.Range("E2").Value = "  -3.97%  "
$ws.Range("D3").Value = "3.335.94"
$ws.Range("E3").Value = "  -1.47%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "573.46"
$ws.Range("E5").Value = "  -3.31%  "
$ws.Range("D6").Value = "181.57"
$ws.Range("E6").Value = "  -5.37%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "0.601"
$ws.Range("E8").Value = "  -1.10%  "
$ws.Range("E9").Value = "  -3.81%  "
$ws.Range("D10").Value = "6.66"
$ws.Range("E10").Value = "  -1.78%  "
$ws.Range("E11").Value = "  -4.52%  "
$ws.Range("D12").Value = "3.915.10"
$ws.Range("E12").Value = "  -1.55%  "
$ws.Range("E13").Value = "  -1.77%  "
$ws.Range("D14").Value = "27.11"
$ws.Range("E14").Value = "  -5.64%  "
$ws.Range("D15").Value = "66.837.18"
$ws.Range("E15").Value = "  -3.93%  "
$ws.Range("E16").Value = "  -2.75%  "
$ws.Range("D17").Value = "3.324.96"
$ws.Range("E17").Value = "  -2.20%  "
$ws.Range("D18").Value = "437.67"
$ws.Range("E18").Value = "  -2.82%  "
$ws.Range("D19").Value = "13.58"
$ws.Range("E19").Value = "  -1.85%  "
$ws.Range("D20").Value = "5.68"
$ws.Range("E20").Value = "  -3.04%  "
$ws.Range("D21").Value = "7.60"
$ws.Range("E21").Value = "  -2.99%  "
$ws.Range("D22").Value = "73.76"
$ws.Range("E22").Value = "  -1.87%  "
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").Value = "0.518"
$ws.Range("E24").Value = "  -0.85%  "
$ws.Range("E25").Value = "  -4.46%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").Value = "9.03"
$ws.Range("E27").Value = "  -4.87%  "
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("E29").Value = "  -1.40%  "
$ws.Range("D30").Value = "22.86"
$ws.Range("E30").Value = "  -2.51%  "
$ws.Range("D31").Value = "5.30"
$ws.Range("E31").Value = "  -6.49%  "
$ws.Range("D32").Value = "0.998"
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("E33").Value = "  -2.79%  "
$ws.Range("D34").Value = "1.23"
$ws.Range("E34").Value = "  -4.69%  "
$ws.Range("D35").Value = "162.32"
$ws.Range("E35").Value = "  -1.82%  "
$ws.Range("D36").Value = "1.49"
$ws.Range("E36").Value = "  -4.31%  "
$ws.Range("D37").Value = "27.88"
$ws.Range("E37").Value = "  +1.18%  "
$ws.Range("D38").Value = "1.84"
$ws.Range("E38").Value = "  -5.83%  "
$ws.Range("D39").Value = "2.820.38"
$ws.Range("E39").Value = "  +2.23%  "
$ws.Range("E40").Value = "  -3.30%  "
$ws.Range("D41").Value = "4.43"
$ws.Range("E41").Value = "  -3.79%  "
$ws.Range("D42").Value = "6.22"
$ws.Range("E42").Value = "  -5.32%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").Value = "40.14"
$ws.Range("E43").Value = "  -1.66%  "
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").Value = "0.0670"
$ws.Range("E44").Value = "  -3.20%  "
$ws.Range("D45").Value = "24.55"
$ws.Range("E45").Value = "  -4.31%  "
$ws.Range("D46").Value = "2.36"
$ws.Range("E46").Value = "  -7.43%  "
$ws.Range("D47").Value = "322.38"
$ws.Range("E47").Value = "  -5.93%  "
$ws.Range("D48").Value = "0.0273"
$ws.Range("E48").Value = "  -4.18%  "
$ws.Range("D49").Value = "0.991"
$ws.Range("E49").Value = "  -3.47%  "
$ws.Range("D50").Value = "6.17"
$ws.Range("E50").Value = "  -3.12%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "0.100"
$ws.Range("E51").Value = "  -1.61%  "
